# Applies the cryptos.xlsx crypto-price/volume/ranking update described by the commit diff.
# Column D/E are scraped as plain text in this sheet (e.g. "214.14", "  -1.54%  "),
# so for any replacement text that Excel would otherwise auto-parse as a number
# (losing formatting like trailing zeros, e.g. "0.1000" -> 0.1), we briefly force
# a text number-format before the assignment, then restore the default cell style
# with Style = "Normal" so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.157.56'
$ws.Range('E2').Value = '  -0.30%  '
$ws.Range('D3').Value = '1.623.86'
$ws.Range('E3').Value = '  -1.30%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '214.25'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.45%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.522'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.42%  '
$ws.Range('E8').Value = '  -1.64%  '
$ws.Range('E9').Value = '  -0.20%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '20.32'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +1.30%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0845'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -0.40%  '
$ws.Range('D12').Value = '1.622.16'
$ws.Range('E12').Value = '  -1.28%  '
$ws.Range('E13').Value = '  -0.42%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.543'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.46%  '
$ws.Range('D15').Value = '27.128.30'
$ws.Range('E15').Value = '  -0.41%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '64.56'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -4.22%  '
$ws.Range('D17').Value = '0.0₃0744'
$ws.Range('E17').Value = '  +0.35%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '216.14'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.67%  '
$ws.Range('E19').Value = '  -0.01%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.93'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.25%  '
$ws.Range('E21').Value = '  -0.98%  '
$ws.Range('E22').Value = '  -6.40%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.05'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.92%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '148.15'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.15%  '
$ws.Range('E25').Value = '  -0.11%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.28'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -3.25%  '
$ws.Range('E27').Value = '  -1.11%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '15.61'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -1.19%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.0506'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.94%  '
$ws.Range('E31').Value = '  -0.57%  '
$ws.Range('E32').Value = '  -1.31%  '
$ws.Range('D33').Value = '1.341.64'
$ws.Range('E33').Value = '  +5.18%  '
$ws.Range('E34').Value = '  -0.58%  '
$ws.Range('E35').Value = '  -0.66%  '
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('E37').Value = '  +1.57%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.860'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.31%  '
$ws.Range('E39').Value = '  -0.08%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.803'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -0.89%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '65.69'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +6.16%  '
$ws.Range('E42').Value = '  -0.35%  '
$ws.Range('E43').Value = '  -1.44%  '
$ws.Range('B44').Value = 'WEMIXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.929'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +38.67%  '
$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D45').Value = '1.760.76'
$ws.Range('E45').Value = '  -1.40%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '90.26'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.81%  '
$ws.Range('E47').Value = '  +0.89%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0105'
$ws.Range('E48').Value = '  -1.68%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.1000'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +2.46%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0513'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.52%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '7.57'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.28%  '
